# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets.
# Both sheets list the same events but the "全部类型" sheet's F23 value
# ends up one higher than the "展览" sheet's (1123 vs 1124).

$wb = $excel.ActiveWorkbook

$commonUpdates = @{
    5  = 252
    9  = 15
    15 = 446
    18 = 404
    19 = 138
    24 = 2844
    27 = 540
    29 = 1610
    32 = 264
    35 = 603
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $commonUpdates[$row]
    }

    if ($sheetName -eq "展览") {
        $ws.Cells.Item(23, 6).Value = 1123
    } else {
        $ws.Cells.Item(23, 6).Value = 1124
    }
}
